$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "email " cell (B1) to drop the trailing space -> "email"
$ws.Range("B1").Value = "email"

# Update the selected cell to match the saved selection state in the file
$ws.Range("F14").Select()
